# Auto-generated Excel COM-interop script
# Applies updated market price / profit cache values to leve profit sheets
# (scheduled runner refresh), per the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 467.78787
$ws.Range("I38").Value = 111.117645
$ws.Range("J38").Value = 846.75
$ws.Range("K38").Value = 333.352935
$ws.Range("L38").Value = 2540.25
$ws.Range("M38").Value = 38.647065
$ws.Range("N38").Value = -3284.25
$ws.Range("H98").Value = 32218.406
$ws.Range("I98").Value = 39405.96
$ws.Range("J98").Value = 1072.3334
$ws.Range("K98").Value = 39405.96
$ws.Range("L98").Value = 1072.3334
$ws.Range("M98").Value = -37907.96
$ws.Range("N98").Value = -4068.3334
$ws.Range("H107").Value = 452.42105
$ws.Range("I107").Value = 388.64706
$ws.Range("J107").Value = 994.5
$ws.Range("K107").Value = 388.64706
$ws.Range("L107").Value = 994.5
$ws.Range("M107").Value = 1531.35294
$ws.Range("N107").Value = -4834.5
$ws.Range("H122").Value = 32218.406
$ws.Range("I122").Value = 39405.96
$ws.Range("J122").Value = 1072.3334
$ws.Range("K122").Value = 118217.88
$ws.Range("L122").Value = 3217.0002
$ws.Range("M122").Value = -115767.88
$ws.Range("N122").Value = -8117.0002
$ws.Range("H129").Value = 41203.24
$ws.Range("I129").Value = 754.8570999999999
$ws.Range("K129").Value = 2264.5713
$ws.Range("M129").Value = 2735.4287
$ws.Range("H132").Value = 1588778.2
$ws.Range("I132").Value = 1852943
$ws.Range("J132").Value = 3789.2
$ws.Range("K132").Value = 5558829
$ws.Range("L132").Value = 11367.6
$ws.Range("M132").Value = -5556299
$ws.Range("N132").Value = -16427.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2207
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2207
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2207
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3019
$ws.Range("H91").Value = 2207
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2207
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2207
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5015
$ws.Range("H122").Value = 1589.9259
$ws.Range("I122").Value = 1535.6923
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4607.0769
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2157.0769
$ws.Range("N122").Value = -13900

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 444.64285
$ws.Range("I80").Value = 69
$ws.Range("K80").Value = 69
$ws.Range("M80").Value = 929
$ws.Range("H83").Value = 444.64285
$ws.Range("I83").Value = 69
$ws.Range("K83").Value = 345
$ws.Range("M83").Value = 4647
$ws.Range("H86").Value = 1554.3846
$ws.Range("I86").Value = 1550.7
$ws.Range("J86").Value = 1566.6666
$ws.Range("K86").Value = 1550.7
$ws.Range("L86").Value = 1566.6666
$ws.Range("M86").Value = -427.7
$ws.Range("N86").Value = -3812.6666
$ws.Range("H89").Value = 1554.3846
$ws.Range("I89").Value = 1550.7
$ws.Range("J89").Value = 1566.6666
$ws.Range("K89").Value = 7753.5
$ws.Range("L89").Value = 7833.333000000001
$ws.Range("M89").Value = -2137.5
$ws.Range("N89").Value = -19065.333
$ws.Range("H94").Value = 750
$ws.Range("I94").Value = 750
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 750
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -299
$ws.Range("N94").ClearContents()
$ws.Range("H107").Value = 1167.2174
$ws.Range("I107").Value = 1012.17645
$ws.Range("J107").Value = 1606.5
$ws.Range("K107").Value = 1012.17645
$ws.Range("L107").Value = 1606.5
$ws.Range("M107").Value = 907.82355
$ws.Range("N107").Value = -5446.5
$ws.Range("H126").Value = 50768
$ws.Range("J126").Value = 50768
$ws.Range("L126").Value = 50768
$ws.Range("N126").Value = -60648
$ws.Range("H134").Value = 1706.5122
$ws.Range("I134").Value = 1569
$ws.Range("K134").Value = 4707
$ws.Range("M134").Value = -2172

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 20898472
$ws.Range("I62").Value = 33435434
$ws.Range("J62").Value = 3533.3333
$ws.Range("K62").Value = 33435434
$ws.Range("L62").Value = 3533.3333
$ws.Range("M62").Value = -33434810
$ws.Range("N62").Value = -4781.3333
$ws.Range("H65").Value = 20898472
$ws.Range("I65").Value = 33435434
$ws.Range("J65").Value = 3533.3333
$ws.Range("K65").Value = 167177170
$ws.Range("L65").Value = 17666.6665
$ws.Range("M65").Value = -167174050
$ws.Range("N65").Value = -23906.6665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 10333.333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 10333.333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 30999.999
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -31539.999
$ws.Range("H67").Value = 10333.333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 10333.333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 30999.999
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -32871.999
$ws.Range("H114").Value = 954.6429000000001
$ws.Range("I114").Value = 290
$ws.Range("J114").Value = 1530.6666
$ws.Range("K114").Value = 870
$ws.Range("L114").Value = 4591.9998
$ws.Range("M114").Value = 2384
$ws.Range("N114").Value = -11099.9998
$ws.Range("H126").Value = 2976
$ws.Range("I126").Value = 1610
$ws.Range("J126").Value = 4342
$ws.Range("K126").Value = 4830
$ws.Range("L126").Value = 13026
$ws.Range("M126").Value = 110
$ws.Range("N126").Value = -22906

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2367.8918
$ws.Range("I102").Value = 1793
$ws.Range("J102").Value = 3726.7273
$ws.Range("K102").Value = 1793
$ws.Range("L102").Value = 3726.7273
$ws.Range("M102").Value = -171
$ws.Range("N102").Value = -6970.7273
$ws.Range("H126").Value = 1920.4166
$ws.Range("I126").Value = 1689
$ws.Range("J126").Value = 2383.25
$ws.Range("K126").Value = 5067
$ws.Range("L126").Value = 7149.75
$ws.Range("M126").Value = -2597
$ws.Range("N126").Value = -12089.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1920.5862
$ws.Range("I7").Value = 1919.88
$ws.Range("J7").Value = 1925
$ws.Range("K7").Value = 1919.88
$ws.Range("L7").Value = 1925
$ws.Range("M7").Value = -1807.88
$ws.Range("N7").Value = -2149
$ws.Range("H40").Value = 2666.5
$ws.Range("I40").Value = 2450
$ws.Range("J40").Value = 3099.5
$ws.Range("K40").Value = 2450
$ws.Range("L40").Value = 3099.5
$ws.Range("M40").Value = -2314
$ws.Range("N40").Value = -3371.5
$ws.Range("H126").Value = 1920.5862
$ws.Range("I126").Value = 1919.88
$ws.Range("J126").Value = 1925
$ws.Range("K126").Value = 5759.64
$ws.Range("L126").Value = 5775
$ws.Range("M126").Value = -3289.64
$ws.Range("N126").Value = -10715

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1692.5714
$ws.Range("I96").Value = 885
$ws.Range("J96").Value = 2769.3333
$ws.Range("K96").Value = 885
$ws.Range("L96").Value = 2769.3333
$ws.Range("M96").Value = 488
$ws.Range("N96").Value = -5515.3333
$ws.Range("H126").Value = 916
$ws.Range("I126").Value = 654.4545000000001
$ws.Range("K126").Value = 1963.3635
$ws.Range("M126").Value = 506.6364999999998
